$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to be stored as text,
# so Excel does not reinterpret numeric-looking strings (with non-standard
# grouping dots, leading/trailing zeros, etc.) as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "71.807.91"
$ws.Range("E2").Value = "  +0.66%  "

# Row 3
$ws.Range("D3").Value = "3.992.99"
$ws.Range("E3").Value = "  -0.41%  "

# Row 4
$ws.Range("E4").Value = "  +0.15%  "

# Row 5
$ws.Range("D5").Value = "528.59"
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("D6").Value = "150.70"
$ws.Range("E6").Value = "  +0.86%  "

# Row 7
$ws.Range("D7").Value = "0.692"
$ws.Range("E7").Value = "  +11.05%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").Value = "0.742"
$ws.Range("E9").Value = "  +0.76%  "

# Row 10
$ws.Range("E10").Value = "  -3.32%  "

# Row 11
$ws.Range("D11").Value = "0.0000323"
$ws.Range("E11").Value = "  -6.10%  "

# Row 12
$ws.Range("D12").Value = "47.22"
$ws.Range("E12").Value = "  +6.58%  "

# Row 13
$ws.Range("D13").Value = "10.54"
$ws.Range("E13").Value = "  -1.13%  "

# Row 14
$ws.Range("D14").Value = "4.627.04"
$ws.Range("E14").Value = "  -0.43%  "

# Row 15
$ws.Range("D15").Value = "3.989.68"
$ws.Range("E15").Value = "  -0.75%  "

# Row 16
$ws.Range("D16").Value = "13.85"
$ws.Range("E16").Value = "  -3.15%  "

# Row 17
$ws.Range("D17").Value = "20.37"
$ws.Range("E17").Value = "  -5.05%  "

# Row 18
$ws.Range("E18").Value = "  -1.24%  "

# Row 19
$ws.Range("E19").Value = "  -3.77%  "

# Row 20
$ws.Range("D20").Value = "71.694.72"
$ws.Range("E20").Value = "  +0.66%  "

# Row 21
$ws.Range("D21").Value = "423.93"
$ws.Range("E21").Value = "  -3.87%  "

# Row 22
$ws.Range("D22").Value = "97.56"
$ws.Range("E22").Value = "  +4.28%  "

# Row 23
$ws.Range("D23").Value = "3.48"
$ws.Range("E23").Value = "  -2.77%  "

# Row 24
$ws.Range("D24").Value = "4.16"
$ws.Range("E24").Value = "  +0.90%  "

# Row 25
$ws.Range("D25").Value = "14.18"
$ws.Range("E25").Value = "  -1.52%  "

# Row 26
$ws.Range("D26").Value = "11.06"
$ws.Range("E26").Value = "  -10.51%  "

# Row 27
$ws.Range("D27").Value = "10.58"
$ws.Range("E27").Value = "  -2.76%  "

# Row 28
$ws.Range("D28").Value = "5.85"
$ws.Range("E28").Value = "  +1.50%  "

# Row 29
$ws.Range("D29").Value = "36.49"
$ws.Range("E29").Value = "  -1.42%  "

# Row 30
$ws.Range("D30").Value = "3.58"
$ws.Range("E30").Value = "  +23.55%  "

# Row 31
$ws.Range("D31").Value = "13.25"
$ws.Range("E31").Value = "  -2.94%  "

# Row 32
$ws.Range("D32").Value = "673.06"
$ws.Range("E32").Value = "  -3.89%  "

# Row 33
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  -1.31%  "

# Row 34
$ws.Range("D34").Value = "6.83"
$ws.Range("E34").Value = "  -0.53%  "

# Row 35
$ws.Range("D35").Value = "65.13"
$ws.Range("E35").Value = "  -2.82%  "

# Row 36
$ws.Range("D36").Value = "41.98"
$ws.Range("E36").Value = "  +1.92%  "

# Row 37
$ws.Range("D37").Value = "0.423"
$ws.Range("E37").Value = "  -4.40%  "

# Row 38
$ws.Range("D38").Value = "0.₃0827"
$ws.Range("E38").Value = "  -9.20%  "

# Row 39
$ws.Range("E39").Value = "  -0.39%  "

# Row 40
$ws.Range("D40").Value = "3.44"
$ws.Range("E40").Value = "  -3.02%  "

# Row 41
$ws.Range("E41").Value = "  -0.27%  "

# Row 42
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.26%  "

# Row 43
$ws.Range("D43").Value = "3.24"
$ws.Range("E43").Value = "  +3.74%  "

# Row 44
$ws.Range("D44").Value = "0.0482"
$ws.Range("E44").Value = "  -2.63%  "

# Row 45
$ws.Range("D45").Value = "0.150"
$ws.Range("E45").Value = "  +2.98%  "

# Row 46
$ws.Range("D46").Value = "9.57"
$ws.Range("E46").Value = "  +2.99%  "

# Row 47
$ws.Range("D47").Value = "2.57"
$ws.Range("E47").Value = "  -12.17%  "

# Row 48
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  -5.32%  "

# Row 49
$ws.Range("D49").Value = "2.97"
$ws.Range("E49").Value = "  -8.53%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "143.67"
$ws.Range("E50").Value = "  +0.22%  "

# Row 51
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "0.000265"
$ws.Range("E51").Value = "  -7.20%  "
